# Applies the commit's changes to the workbook:
#  1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" worksheet (sheetId 3) with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1 & 2: header renames (cells already carry the bold/centered/bordered style) ---
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, centered, bordered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$newSheet.Cells.Item(2,1).Value = 45501.99999999999
$newSheet.Cells.Item(2,2).Value = 87
$newSheet.Cells.Item(2,3).Value = -133.5700279743517
$newSheet.Cells.Item(2,4).Value = 317.3979416077502
$newSheet.Cells.Item(3,1).Value = 45508.99999999999
$newSheet.Cells.Item(3,2).Value = 107
$newSheet.Cells.Item(3,3).Value = -123.8732781389888
$newSheet.Cells.Item(3,4).Value = 354.1347296291844
$newSheet.Cells.Item(4,1).Value = 45515.99999999999
$newSheet.Cells.Item(4,2).Value = 127
$newSheet.Cells.Item(4,3).Value = -98.84325266268957
$newSheet.Cells.Item(4,4).Value = 330.6341736856485
$newSheet.Cells.Item(5,1).Value = 45522.99999999999
$newSheet.Cells.Item(5,2).Value = 147
$newSheet.Cells.Item(5,3).Value = -88.23342132810231
$newSheet.Cells.Item(5,4).Value = 358.8846980453137
$newSheet.Cells.Item(6,1).Value = 45529.99999999999
$newSheet.Cells.Item(6,2).Value = 168
$newSheet.Cells.Item(6,3).Value = -57.94393930528148
$newSheet.Cells.Item(6,4).Value = 384.832129021496
$newSheet.Cells.Item(7,1).Value = 45536.99999999999
$newSheet.Cells.Item(7,2).Value = 188
$newSheet.Cells.Item(7,3).Value = -56.31349919330584
$newSheet.Cells.Item(7,4).Value = 413.6689129929982
$newSheet.Cells.Item(8,1).Value = 45543.99999999999
$newSheet.Cells.Item(8,2).Value = 208
$newSheet.Cells.Item(8,3).Value = -1.993894224750538
$newSheet.Cells.Item(8,4).Value = 441.2317411006823
$newSheet.Cells.Item(9,1).Value = 45550.99999999999
$newSheet.Cells.Item(9,2).Value = 228
$newSheet.Cells.Item(9,3).Value = -1.330933838025009
$newSheet.Cells.Item(9,4).Value = 444.6385717943772
$newSheet.Cells.Item(10,1).Value = 45557.99999999999
$newSheet.Cells.Item(10,2).Value = 248
$newSheet.Cells.Item(10,3).Value = 34.25300074221975
$newSheet.Cells.Item(10,4).Value = 480.4806326040288
$newSheet.Cells.Item(11,1).Value = 45585.99999999999
$newSheet.Cells.Item(11,2).Value = 329
$newSheet.Cells.Item(11,3).Value = 102.1399993988326
$newSheet.Cells.Item(11,4).Value = 533.2430215000944
$newSheet.Cells.Item(12,1).Value = 45599.99999999999
$newSheet.Cells.Item(12,2).Value = 369
$newSheet.Cells.Item(12,3).Value = 143.3458093666607
$newSheet.Cells.Item(12,4).Value = 574.7170930668999
$newSheet.Cells.Item(13,1).Value = 45606.99999999999
$newSheet.Cells.Item(13,2).Value = 389
$newSheet.Cells.Item(13,3).Value = 165.803955793214
$newSheet.Cells.Item(13,4).Value = 613.3825905087509
$newSheet.Cells.Item(14,1).Value = 45613.99999999999
$newSheet.Cells.Item(14,2).Value = 409
$newSheet.Cells.Item(14,3).Value = 184.9450566737841
$newSheet.Cells.Item(14,4).Value = 635.0451747205949
$newSheet.Cells.Item(15,1).Value = 45620.99999999999
$newSheet.Cells.Item(15,2).Value = 429
$newSheet.Cells.Item(15,3).Value = 210.6434867499835
$newSheet.Cells.Item(15,4).Value = 642.6080624430505
$newSheet.Cells.Item(16,1).Value = 45627.99999999999
$newSheet.Cells.Item(16,2).Value = 450
$newSheet.Cells.Item(16,3).Value = 226.6576978826955
$newSheet.Cells.Item(16,4).Value = 680.8102832295194
$newSheet.Cells.Item(17,1).Value = 45634.99999999999
$newSheet.Cells.Item(17,2).Value = 470
$newSheet.Cells.Item(17,3).Value = 230.0480007900755
$newSheet.Cells.Item(17,4).Value = 690.5837396814924
$newSheet.Cells.Item(18,1).Value = 45641.99999999999
$newSheet.Cells.Item(18,2).Value = 490
$newSheet.Cells.Item(18,3).Value = 258.2105187049225
$newSheet.Cells.Item(18,4).Value = 713.8298802014922
$newSheet.Cells.Item(19,1).Value = 45648.99999999999
$newSheet.Cells.Item(19,2).Value = 510
$newSheet.Cells.Item(19,3).Value = 295.6340102407407
$newSheet.Cells.Item(19,4).Value = 726.0535264217915
$newSheet.Cells.Item(20,1).Value = 45655.99999999999
$newSheet.Cells.Item(20,2).Value = 530
$newSheet.Cells.Item(20,3).Value = 298.6876042698062
$newSheet.Cells.Item(20,4).Value = 755.3116167262162

# Apply the date style (as used on other sheets' date column) to column A of the data rows
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A20").PasteSpecial(-4122)

$newSheet.Range("A1").Select()

Write-Host "PO Forecast sheet added and headers renamed."
